# New crime data collected - refresh of the weekly CompStat numbers
# (20th Precinct report, week covering 7/31/2023 - 8/6/2023,
# Volume 30 Number 31).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header (masthead) text ----------------------------------------
$ws.Range('A8').Value2 = 'Volume 30   Number  31'
$ws.Range('C9').Value2 = 'Report Covering the Week  7/31/2023  Through  8/6/2023'

# --- Crime-complaint table (rows 15-27) -------------------------------------
# Row 15: Rape
$ws.Range('N15').Value2 = -64.705882352941

# Row 16: Robbery
$ws.Range('C16').Value2 = 3
$ws.Range('D16').Value2 = 5
$ws.Range('E16').Value2 = -40
$ws.Range('F16').Value2 = 6
$ws.Range('G16').Value2 = 12
$ws.Range('H16').Value2 = -50
$ws.Range('I16').Value2 = 57
$ws.Range('J16').Value2 = 74
$ws.Range('K16').Value2 = -22.972972972973
$ws.Range('L16').Value2 = -3.389830508474
$ws.Range('M16').Value2 = 1.785714285714
$ws.Range('N16').Value2 = -84.383561643835

# Row 17: Fel. Assault
$ws.Range('F17').Value2 = 8
$ws.Range('H17').Value2 = 100
$ws.Range('I17').Value2 = 67
$ws.Range('J17').Value2 = 59
$ws.Range('K17').Value2 = 13.559322033898
$ws.Range('L17').Value2 = 97.058823529411
$ws.Range('M17').Value2 = 52.272727272727
$ws.Range('N17').Value2 = 9.83606557377

# Row 18: Burglary
$ws.Range('C18').Value2 = 3
$ws.Range('D18').Value2 = 6
$ws.Range('F18').Value2 = 7
$ws.Range('G18').Value2 = 19
$ws.Range('H18').Value2 = -63.157894736842
$ws.Range('I18').Value2 = 70
$ws.Range('J18').Value2 = 98
$ws.Range('K18').Value2 = -28.571428571428
$ws.Range('L18').Value2 = 27.272727272727
$ws.Range('M18').Value2 = 4.477611940298
$ws.Range('N18').Value2 = -86.590038314176

# Row 19: Gr. Larceny
$ws.Range('C19').Value2 = 14
$ws.Range('D19').Value2 = 8
$ws.Range('E19').Value2 = 75
$ws.Range('F19').Value2 = 58
$ws.Range('G19').Value2 = 67
$ws.Range('H19').Value2 = -13.432835820895
$ws.Range('I19').Value2 = 430
$ws.Range('J19').Value2 = 480
$ws.Range('K19').Value2 = -10.416666666666
$ws.Range('L19').Value2 = 15.281501340482
$ws.Range('M19').Value2 = 10.824742268041
$ws.Range('N19').Value2 = -62.770562770562

# Row 20: G.L.A.
$ws.Range('D20').Value2 = 2
$ws.Range('E20').Value2 = -50
$ws.Range('G20').Value2 = 11
$ws.Range('H20').Value2 = 36.363636363636
$ws.Range('I20').Value2 = 65
$ws.Range('J20').Value2 = 40
$ws.Range('K20').Value2 = 62.5
$ws.Range('L20').Value2 = 18.181818181818
$ws.Range('M20').Value2 = 225
$ws.Range('N20').Value2 = -90.412979351032

# Row 21: TOTAL
$ws.Range('C21').Value2 = 22
$ws.Range('D21').Value2 = 22
$ws.Range('E21').Value2 = 0
$ws.Range('G21').Value2 = 116
$ws.Range('H21').Value2 = -18.103448275862
$ws.Range('I21').Value2 = 696
$ws.Range('J21').Value2 = 761
$ws.Range('K21').Value2 = -8.541392904073
$ws.Range('L21').Value2 = 19.793459552495
$ws.Range('M21').Value2 = 19.382504288164
$ws.Range('N21').Value2 = -75.13397642015

# Row 22: Transit
$ws.Range('J22').Value2 = 19
$ws.Range('K22').Value2 = -15.78947368421

# Row 23: Housing
$ws.Range('F23').Value2 = 1
$ws.Range('G23').Value2 = 2
$ws.Range('H23').Value2 = -50
$ws.Range('J23').Value2 = 15
$ws.Range('K23').Value2 = 66.666666666666

# Row 24: Petit Larceny
$ws.Range('C24').Value2 = 19
$ws.Range('D24').Value2 = 19
$ws.Range('E24').Value2 = 0
$ws.Range('G24').Value2 = 88
$ws.Range('H24').Value2 = 26.136363636363
$ws.Range('I24').Value2 = 656
$ws.Range('J24').Value2 = 800
$ws.Range('K24').Value2 = -18
$ws.Range('L24').Value2 = -20.581113801452
$ws.Range('M24').Value2 = 7.18954248366

# Row 25: Misd. Assault
$ws.Range('C25').Value2 = 3
$ws.Range('D25').Value2 = 2
$ws.Range('E25').Value2 = 50
$ws.Range('F25').Value2 = 19
$ws.Range('G25').Value2 = 12
$ws.Range('H25').Value2 = 58.333333333333
$ws.Range('I25').Value2 = 141
$ws.Range('J25').Value2 = 133
$ws.Range('K25').Value2 = 6.015037593984
$ws.Range('L25').Value2 = 39.603960396039
$ws.Range('M25').Value2 = -17.058823529411

# Row 26: UCR Rape*
$ws.Range('G26').Value2 = 4
$ws.Range('H26').Value2 = -75

# Row 27: Other Sex Crimes
$ws.Range('D27').Value2 = 1
$ws.Range('D27').NumberFormat = "#,##0"
$ws.Range('E27').Value2 = 0
$ws.Range('E27').NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range('G27').Value2 = 2
$ws.Range('H27').Value2 = 100
$ws.Range('I27').Value2 = 25
$ws.Range('J27').Value2 = 37
$ws.Range('K27').Value2 = -32.432432432432
$ws.Range('L27').Value2 = -7.407407407407
